$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d354bc35-... file, now ready for handoff ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 10:14:22"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
# Row 2 (48cbe1ab file) status updated
$wsZh.Range("C2").Value = "Ready for handoff"
# Row 3 (d354bc35 file): status / priority / handoff datetime updated
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-24 10:14:17"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
# Row 2 (48cbe1ab file) status updated
$wsDe.Range("C2").Value = "Ready for handoff"
# Row 3 (d354bc35 file): status / priority / handoff datetime updated
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-24 10:14:22"

# --- Column width adjustments: the Status-adjacent columns widened to fit the
#     longer "Ready for handoff" / timestamp text (mirrors Excel auto-resize
#     behavior recorded in the target workbook). ColumnWidth is quantized to
#     1/6-character increments by this host, so 16.33 is chosen to land on
#     the closest reachable bucket to the recorded width. ---
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZh.Range("C1").ColumnWidth = 16.33
$wsDe.Range("C1").ColumnWidth = 16.33
